$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as plain text, preserving number-like strings
# (e.g. "318.51", "1.00", "0.0920") exactly instead of letting Excel coerce
# them into numeric cells.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '45.478.77'
$ws.Range("E2").Value = '  +0.10%  '

Set-TextValue "D3" '2.378.78'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  -0.07%  '

Set-TextValue "D5" '318.51'
$ws.Range("E5").Value = '  +0.40%  '

Set-TextValue "D6" '109.41'
$ws.Range("E6").Value = '  -3.50%  '

Set-TextValue "D7" '0.639'
$ws.Range("E7").Value = '  +0.94%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  -1.40%  '

Set-TextValue "D10" '41.06'
$ws.Range("E10").Value = '  -3.72%  '

Set-TextValue "D11" '0.0920'
$ws.Range("E11").Value = '  -1.02%  '

Set-TextValue "D12" '8.56'
$ws.Range("E12").Value = '  -1.06%  '

Set-TextValue "D14" '0.988'
$ws.Range("E14").Value = '  -2.85%  '

Set-TextValue "D15" '2.740.90'
$ws.Range("E15").Value = '  -0.04%  '

Set-TextValue "D16" '15.46'
$ws.Range("E16").Value = '  -2.45%  '

Set-TextValue "D17" '2.384.30'
$ws.Range("E17").Value = '  +0.16%  '

Set-TextValue "D18" '45.407.01'
$ws.Range("E18").Value = '  +0.19%  '

Set-TextValue "D19" '15.67'
$ws.Range("E19").Value = '  +17.34%  '

$ws.Range("E20").Value = '  -3.72%  '

$ws.Range("E21").Value = '  -0.96%  '

Set-TextValue "D22" '3.72'
$ws.Range("E22").Value = '  +5.42%  '

Set-TextValue "D23" '73.37'
$ws.Range("E23").Value = '  -1.93%  '

Set-TextValue "D24" '261.56'
$ws.Range("E24").Value = '  -2.87%  '

Set-TextValue "D25" '2.36'
$ws.Range("E25").Value = '  +0.18%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("E27").Value = '  +0.09%  '

Set-TextValue "D28" '11.29'
$ws.Range("E28").Value = '  -0.18%  '

Set-TextValue "D29" '2.29'
$ws.Range("E29").Value = '  -1.68%  '

Set-TextValue "D30" '22.49'
$ws.Range("E30").Value = '  -1.80%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D31" '37.64'
$ws.Range("E31").Value = '  -3.89%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D32" '0.0953'
$ws.Range("E32").Value = '  -1.29%  '

Set-TextValue "D33" '168.43'
$ws.Range("E33").Value = '  -1.99%  '

$ws.Range("E34").Value = '  -2.06%  '

Set-TextValue "D35" '0.133'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("E36").Value = '  -2.30%  '

$ws.Range("E37").Value = '  -3.29%  '

Set-TextValue "D38" '1.96'
$ws.Range("E38").Value = '  +13.79%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D39" '3.03'
$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("B40").Value = 'NEARProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D40" '4.01'
$ws.Range("E40").Value = '  -0.04%  '

Set-TextValue "D41" '0.0357'
$ws.Range("E41").Value = '  -2.28%  '

Set-TextValue "D42" '97.49'
$ws.Range("E42").Value = '  -7.07%  '

Set-TextValue "D43" '70.56'
$ws.Range("E43").Value = '  -1.46%  '

Set-TextValue "D44" '0.229'
$ws.Range("E44").Value = '  -4.01%  '

Set-TextValue "D45" '12.99'
$ws.Range("E45").Value = '  -1.73%  '

Set-TextValue "D46" '1.852.81'
$ws.Range("E46").Value = '  +12.88%  '

$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue "D47" '6.00'
$ws.Range("E47").Value = '  +3.47%  '

$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D48" '1.00'
$ws.Range("E48").Value = '  +0.01%  '

Set-TextValue "D49" '83.84'
$ws.Range("E49").Value = '  +5.95%  '

Set-TextValue "D50" '112.87'
$ws.Range("E50").Value = '  -3.12%  '

Set-TextValue "D51" '9.28'
$ws.Range("E51").Value = '  -0.66%  '
